# Test data update for 4.0 release
# Rename the "Interventional" study-type expected-filename entries to
# "Clinical" on Sheet1 (rows 3 and 4, column G).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G3").Value = "ExcelReport-NewImportLogic_1 - Test_Automation_1-Clinical-"
$ws.Range("G4").Value = "WordReport-NewImportLogic_1 - Test_Automation_1-Clinical-"

$ws.Range("G4").Select()
